$d = $word.ActiveDocument

# --- helper: simple exact-text Find & Replace (no wildcards) ---
function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) Add "US Citizen   ●   " before the "Bay Area ..." contact line.
$rng = $d.Content
$rng.Find.Execute("Bay Area   ") | Out-Null
$insertPoint = $rng.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertBefore("US Citizen   ●   ")

# 2) Update the SKILLS line: add Postgres & Express, drop Azure, restructure Agile(Scrum).
Replace-Text "HTML, CSS, SQL, MongoDB, Node.js, Flask, Git, Linux, AWS, Azure, Agile (Scrum)" `
             "HTML, CSS, SQL, Postgres, MongoDB, Node.js, Express, Flask, Git, Linux, AWS, Scrum, Agile"

# 3 & 5) "Software Engineering Intern" -> "Software Engineer Intern" (both occurrences).
Replace-Text "Software Engineering Intern" "Software Engineer Intern"

# 4) Remove the stray _GoBack bookmark left on the Coursera bullet (it gets relocated below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 6) Project ANT bullet rewrite.
Replace-Text "Developed back-end architecture for Project ANT, a freelancing platform for building and showcasing user portfolios. " `
             "Lead the portfolio team for Project ANT, a freelancing platform for building and showcasing user portfolios. "

# 7) Databases bullet rewrite.
Replace-Text "Designed and maintained NoSQL and MongoDB databases and improved web application with Node.js and Express. " `
             "Designed and implemented Postgres databases for digital portfolios using Elastic Beanstalk, Node.js, and Express. "

# 8) "Volunteer Developer  -  " -> "Software Engineer Intern  -  "
Replace-Text "Volunteer Developer  -  " "Software Engineer Intern  -  "

# 9) "CARE International" (the org/title run) -> "Develop for Good"
Replace-Text "CARE International" "Develop for Good"

# 10) Re-add the _GoBack bookmark right after the tab, before the "September 2020 - Present" date
#     on the (now) "Develop for Good" bullet line (the 2nd "September 2020 - Present" in the doc).
$titleRng = $d.Content
$titleRng.Find.Execute("Develop for Good") | Out-Null
$afterTitle = $d.Range($titleRng.End, $d.Content.End)
$afterTitle.Find.Execute("September 2020") | Out-Null
$bmPoint = $afterTitle.Duplicate
$bmPoint.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# 11) CARE Yemen bullet rewrite.
Replace-Text "Launched interactive mobile application for the CARE Yemen Youth and Women Initiative to collect and transfer data from local authorities to youth and women initiatives and service providers. " `
             "Launched interactive mobile application for the CARE International Yemen Youth and Women Initiative to collect and transfer data from local authorities to youth and women initiatives, service providers, and other target groups."

# 12) Hackathon placement bullet rewrite.
Replace-Text "Third place at the Open Water VC Fall Hackathon." "Third place at the Open Water Accelerator Fall Hackathon."

# 13) Volunteering bullet rewrite.
Replace-Text "Volunteered for future hackathons as staff to mentor and technically support all participants." `
             "Volunteered for future hackathons as staff to mentor and technically support over one hundred participants. "
